$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.299.19"
$ws.Range("E2").Value = "  +0.42%  "

$ws.Range("D3").Value = "2.014.59"
$ws.Range("E3").Value = "  -1.03%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.87"
$ws.Range("E5").Value = "  +3.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.642"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.06"
$ws.Range("E7").Value = "  +12.27%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.28"
$ws.Range("E9").Value = "  -4.93%  "

$ws.Range("E10").Value = "  +1.71%  "

$ws.Range("E11").Value = "  +0.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.105"
$ws.Range("E12").Value = "  -1.50%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.909"
$ws.Range("E13").Value = "  +1.09%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.84"
$ws.Range("E14").Value = "  +5.11%  "

$ws.Range("D15").Value = "2.306.71"
$ws.Range("E15").Value = "  -1.05%  "

$ws.Range("E16").Value = "  +1.55%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.86"
$ws.Range("E17").Value = "  +14.49%  "

$ws.Range("D18").Value = "2.006.09"
$ws.Range("E18").Value = "  -1.41%  "

$ws.Range("D19").Value = "36.239.47"
$ws.Range("E19").Value = "  +0.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.12"
$ws.Range("E20").Value = "  +1.28%  "

$ws.Range("E21").Value = "  +1.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.28"
$ws.Range("E22").Value = "  +2.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.31"
$ws.Range("E23").Value = "  -1.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.66"
$ws.Range("E24").Value = "  +18.61%  "

$ws.Range("E25").Value = "  +0.16%  "

$ws.Range("E26").Value = "  -1.76%  "

$ws.Range("E27").Value = "  +4.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.05"
$ws.Range("E28").Value = "  -0.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.64"
$ws.Range("E29").Value = "  -1.11%  "

$ws.Range("E30").Value = "  +35.01%  "

$ws.Range("E31").Value = "  -0.06%  "

$ws.Range("E32").Value = "  +3.62%  "

$ws.Range("E33").Value = "  -0.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.56"
$ws.Range("E34").Value = "  +3.92%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0608"
$ws.Range("E35").Value = "  +2.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.48"
$ws.Range("E36").Value = "  +13.41%  "

$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("E38").Value = "  -0.61%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.90"
$ws.Range("E39").Value = "  +17.72%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.103"
$ws.Range("E40").Value = "  +14.66%  "

$ws.Range("E41").Value = "  +0.94%  "

$ws.Range("E42").Value = "  +1.23%  "

$ws.Range("E43").Value = "  +1.46%  "

$ws.Range("E44").Value = "  +3.49%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.70"
$ws.Range("E45").Value = "  +6.13%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.88"
$ws.Range("E46").Value = "  +7.90%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "94.36"
$ws.Range("E47").Value = "  +1.77%  "

$ws.Range("D48").Value = "1.431.22"
$ws.Range("E48").Value = "  +5.56%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.59"
$ws.Range("E49").Value = "  +15.92%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.92"
$ws.Range("E50").Value = "  -0.32%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.39"
$ws.Range("E51").Value = "  +4.38%  "
